$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 56353
$ws.Range("D2").Value = 115169332
$ws.Range("C3").Value = 136665
$ws.Range("D3").Value = 321483795
$ws.Range("C4").Value = 49627
$ws.Range("D4").Value = 144421780
$ws.Range("C5").Value = 15719
$ws.Range("D5").Value = 53279480
$ws.Range("C6").Value = 5793
$ws.Range("D6").Value = 26423879
$ws.Range("C7").Value = 1159
$ws.Range("D7").Value = 6748100
$ws.Range("C8").Value = 72
$ws.Range("D8").Value = 491481
$ws.Range("C12").Value = 58679
$ws.Range("D12").Value = 93468373
$ws.Range("C13").Value = 14298
$ws.Range("D13").Value = 28882200
$ws.Range("C14").Value = 38452
$ws.Range("D14").Value = 88882866
$ws.Range("C15").Value = 12773
$ws.Range("D15").Value = 35414928
$ws.Range("C16").Value = 3681
$ws.Range("D16").Value = 11326601
$ws.Range("C17").Value = 1204
$ws.Range("D17").Value = 5198082
$ws.Range("C18").Value = 269
$ws.Range("D18").Value = 1454775
$ws.Range("C19").Value = 15
$ws.Range("D19").Value = 72134
$ws.Range("C20").Value = 14354
$ws.Range("D20").Value = 22262188
$ws.Range("C21").Value = 20087
$ws.Range("D21").Value = 42566121
$ws.Range("C22").Value = 47582
$ws.Range("D22").Value = 114554592
$ws.Range("C23").Value = 16475
$ws.Range("D23").Value = 47563271
$ws.Range("C24").Value = 4878
$ws.Range("D24").Value = 15791129
$ws.Range("C25").Value = 1574
$ws.Range("D25").Value = 6598224
$ws.Range("C26").Value = 250
$ws.Range("D26").Value = 1371338
$ws.Range("C28").Value = 16070
$ws.Range("D28").Value = 24785102
$ws.Range("C29").Value = 11365
$ws.Range("D29").Value = 23274259
$ws.Range("C30").Value = 32899
$ws.Range("D30").Value = 74473962
$ws.Range("C31").Value = 11873
$ws.Range("D31").Value = 32205544
$ws.Range("C32").Value = 3271
$ws.Range("D32").Value = 9857112
$ws.Range("C33").Value = 1025
$ws.Range("D33").Value = 4346045
$ws.Range("C34").Value = 207
$ws.Range("D34").Value = 963490
$ws.Range("C35").Value = 8
$ws.Range("D35").Value = 36015
$ws.Range("C36").Value = 11622
$ws.Range("D36").Value = 18057235
$ws.Range("C37").Value = 5034
$ws.Range("D37").Value = 11144416
$ws.Range("C38").Value = 11811
$ws.Range("D38").Value = 27527295
$ws.Range("C39").Value = 4881
$ws.Range("D39").Value = 13923245
$ws.Range("C41").Value = 434
$ws.Range("D41").Value = 2077684
$ws.Range("C42").Value = 55
$ws.Range("D42").Value = 359772
$ws.Range("C44").Value = 3580
$ws.Range("D44").Value = 5518836
$ws.Range("C45").Value = 25766
$ws.Range("D45").Value = 53356992
$ws.Range("C46").Value = 76438
$ws.Range("D46").Value = 180631620
$ws.Range("C47").Value = 29405
$ws.Range("D47").Value = 82216053
$ws.Range("C48").Value = 9575
$ws.Range("D48").Value = 29264517
$ws.Range("C49").Value = 3304
$ws.Range("D49").Value = 13329859
$ws.Range("C50").Value = 579
$ws.Range("D50").Value = 3259376
$ws.Range("C53").Value = 26509
$ws.Range("D53").Value = 48807837
$ws.Range("C54").Value = 2709
$ws.Range("D54").Value = 4390992
$ws.Range("C55").Value = 9042
$ws.Range("D55").Value = 15004234
$ws.Range("C57").Value = 995
$ws.Range("D57").Value = 1951844
$ws.Range("C61").Value = 9259
$ws.Range("D61").Value = 13776525
$ws.Range("C62").Value = 1826
$ws.Range("D62").Value = 3990094
$ws.Range("C63").Value = 4325
$ws.Range("D63").Value = 9408104
$ws.Range("C64").Value = 1729
$ws.Range("D64").Value = 3894584
$ws.Range("C66").Value = 211
$ws.Range("D66").Value = 477383
$ws.Range("C68").Value = 2833
$ws.Range("D68").Value = 5606325
$ws.Range("C69").Value = 22931
$ws.Range("D69").Value = 45434148
$ws.Range("C70").Value = 66545
$ws.Range("D70").Value = 152174983
$ws.Range("C71").Value = 24349
$ws.Range("D71").Value = 67876410
$ws.Range("C72").Value = 7630
$ws.Range("D72").Value = 23218719
$ws.Range("C73").Value = 2472
$ws.Range("D73").Value = 10010043
$ws.Range("C74").Value = 505
$ws.Range("D74").Value = 2813570
$ws.Range("C78").Value = 21273
$ws.Range("D78").Value = 32710617
$ws.Range("C79").Value = 83669
$ws.Range("D79").Value = 172328303
$ws.Range("C80").Value = 227295
$ws.Range("D80").Value = 513898224
$ws.Range("C81").Value = 102696
$ws.Range("D81").Value = 289145226
$ws.Range("C82").Value = 37469
$ws.Range("D82").Value = 126508655
$ws.Range("C83").Value = 13839
$ws.Range("D83").Value = 62716905
$ws.Range("C84").Value = 2719
$ws.Range("D84").Value = 17510662
$ws.Range("C90").Value = 79406
$ws.Range("D90").Value = 126327182
$ws.Range("C91").Value = 5642
$ws.Range("D91").Value = 8808893
$ws.Range("C92").Value = 13582
$ws.Range("D92").Value = 21512438
$ws.Range("C93").Value = 4359
$ws.Range("D93").Value = 7096451
$ws.Range("C98").Value = 6370
$ws.Range("D98").Value = 8768095
$ws.Range("C99").Value = 2311
$ws.Range("D99").Value = 4280299
$ws.Range("C100").Value = 7346
$ws.Range("D100").Value = 14632947
$ws.Range("C101").Value = 2632
$ws.Range("D101").Value = 5972509
$ws.Range("C102").Value = 976
$ws.Range("D102").Value = 2344417
$ws.Range("C103").Value = 340
$ws.Range("D103").Value = 1090276
$ws.Range("C104").Value = 69
$ws.Range("D104").Value = 314891
$ws.Range("C106").Value = 4926
$ws.Range("D106").Value = 7208094
$ws.Range("C107").Value = 1040
$ws.Range("D107").Value = 2364697
$ws.Range("C108").Value = 725
$ws.Range("D108").Value = 1835936
$ws.Range("C113").Value = 16425
$ws.Range("D113").Value = 34555192
$ws.Range("C114").Value = 43249
$ws.Range("D114").Value = 101493644
$ws.Range("C115").Value = 15336
$ws.Range("D115").Value = 42956661
$ws.Range("C116").Value = 4738
$ws.Range("D116").Value = 15021472
$ws.Range("C117").Value = 1473
$ws.Range("D117").Value = 6268648
$ws.Range("C118").Value = 305
$ws.Range("D118").Value = 1735408
$ws.Range("C122").Value = 13544
$ws.Range("D122").Value = 20793975
$ws.Range("C123").Value = 44777
$ws.Range("D123").Value = 90888103
$ws.Range("C124").Value = 95829
$ws.Range("D124").Value = 216140150
$ws.Range("C125").Value = 32396
$ws.Range("D125").Value = 87300513
$ws.Range("C126").Value = 10157
$ws.Range("D126").Value = 30950053
$ws.Range("C127").Value = 3226
$ws.Range("D127").Value = 13308333
$ws.Range("C128").Value = 662
$ws.Range("D128").Value = 3604327
$ws.Range("C129").Value = 33
$ws.Range("D129").Value = 173657
$ws.Range("C132").Value = 35107
$ws.Range("D132").Value = 53741856
$ws.Range("C133").Value = 53878
$ws.Range("D133").Value = 110716034
$ws.Range("C134").Value = 112820
$ws.Range("D134").Value = 252439888
$ws.Range("C135").Value = 36537
$ws.Range("D135").Value = 100766268
$ws.Range("C136").Value = 10819
$ws.Range("D136").Value = 33493378
$ws.Range("C137").Value = 3419
$ws.Range("D137").Value = 14132899
$ws.Range("C138").Value = 560
$ws.Range("D138").Value = 3108229
$ws.Range("C139").Value = 48
$ws.Range("D139").Value = 234414
$ws.Range("C142").Value = 43853
$ws.Range("D142").Value = 65784776
$ws.Range("C143").Value = 19680
$ws.Range("D143").Value = 40476153
$ws.Range("C144").Value = 47839
$ws.Range("D144").Value = 112779715
$ws.Range("C145").Value = 17974
$ws.Range("D145").Value = 50421785
$ws.Range("C146").Value = 5179
$ws.Range("D146").Value = 16038485
$ws.Range("C147").Value = 1552
$ws.Range("D147").Value = 6643390
$ws.Range("C148").Value = 346
$ws.Range("D148").Value = 1995266
$ws.Range("C152").Value = 14827
$ws.Range("D152").Value = 23041689
$ws.Range("C153").Value = 53409
$ws.Range("D153").Value = 110725542
$ws.Range("C154").Value = 124062
$ws.Range("D154").Value = 286850289
$ws.Range("C155").Value = 39478
$ws.Range("D155").Value = 113520481
$ws.Range("C156").Value = 11807
$ws.Range("D156").Value = 39852506
$ws.Range("C157").Value = 4230
$ws.Range("D157").Value = 18954236
$ws.Range("C158").Value = 861
$ws.Range("D158").Value = 5252801
$ws.Range("C159").Value = 53
$ws.Range("D159").Value = 256597
$ws.Range("C160").Value = 41149
$ws.Range("D160").Value = 63998965
